$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Mark all currently un-marked inline drawings as "no proofing needed"
#    (mirrors Word re-flagging pictures with <w:noProof/> on save after
#    the document passes through a proofing/edit pass). The first two
#    inline shapes already carry <w:noProof/>; shapes 3..Count do not.
# ---------------------------------------------------------------------
for ($i = 3; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    $shp.Range.NoProofing = $true
}

# ---------------------------------------------------------------------
# 2) "Click on line NUMBER 22 above to establish a breakpoint."
#       -> "Click on line number 21 above to establish a breakpoint."
#    Done as two separate retypes (the digits, then the word "NUMBER")
#    so the final run layout / _GoBack bookmark placement matches a
#    natural, interactive Word edit.
# ---------------------------------------------------------------------

# 2a) " 22" -> " 21"
$r1 = $d.Content
$r1.Find.ClearFormatting()
$r1.Find.Text = " 22"
$r1.Find.Forward = $true
$r1.Find.Wrap = 0
$r1.Find.Execute() | Out-Null
$r1.Text = " 21"

# 2b) "NUMBER" -> "number" (this is the last edit, so it is where Word
#     drops the _GoBack bookmark)
$r2 = $d.Content
$r2.Find.ClearFormatting()
$r2.Find.Text = "NUMBER"
$r2.Find.Forward = $true
$r2.Find.Wrap = 0
$r2.Find.Execute() | Out-Null
$r2.Text = "number"

# Move the document's _GoBack bookmark to right after "number" (Word
# keeps only a single _GoBack, always at the most recent edit point).
$goBackRange = $d.Range($r2.End, $r2.End)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# 2c) Force " 21" / " above to establish a breakpoint." onto separate
#     runs (matching the two distinct retype actions) by briefly
#     bookmarking the boundary between them.
$r3 = $d.Range($r2.End, $d.Content.End)
$r3.Find.ClearFormatting()
$r3.Find.Text = " 21"
$r3.Find.Forward = $true
$r3.Find.Wrap = 0
$r3.Find.Execute() | Out-Null

$splitRange = $d.Range($r3.End, $r3.End)
$d.Bookmarks.Add("ZZTempSplit", $splitRange)
$d.Bookmarks.Item("ZZTempSplit").Delete()

# ---------------------------------------------------------------------
# 3) Footer "Page X of Y" cached PAGE field result: 1 -> 5 (Word
#    recalculates/recaches this display text whenever the document is
#    resaved).
# ---------------------------------------------------------------------
$ftr = $d.Sections.Item(1).Footers.Item(1)
$pageResultChar = $ftr.Range.Characters.Item(8)
if ($pageResultChar.Text -eq "1") {
    $pageResultChar.Text = "5"
}
